# "Pruebas untarias y registro de tiempos" - update registered time
# estimates (column F) on a few use-case rows and leave the sheet
# scrolled/selected where the author last left it (F30).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")
$ws.Activate()

# Registro de tiempos (hours) for the affected use cases.
$ws.Range("F5").Value = 19.440000000000001
$ws.Range("F6").Value = 2.6600000000000001
$ws.Range("F7").Value = 2.2200000000000002
$ws.Range("F29").Value = 7.1600000000000001

# Leave the view scrolled to where row 21 is at the top and F30 selected,
# matching where the author was working when the file was saved.
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 21
$ws.Range("F30").Select()
